$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "SwordMeta"

$ws.Range("D1").Value = "最大间隔"
$ws.Range("E1").Value = "描述"
$ws.Range("F1").Value = "动画索引"

$ws.Range("D2").Value = "maxInterval"
$ws.Range("E2").Value = "description"
$ws.Range("F2").Value = "animIDs"

$ws.Range("D3").Value = "float"
$ws.Range("E3").Value = "string"
$ws.Range("F3").Value = "int[]"

$ws.Range("B4").Value = "绿色的剑"
$ws.Range("D4").Value = 1.5
$ws.Range("E4").Value = "绿色的一把剑"
$ws.Range("F4").Value = "0,1,2"

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "红色的剑"
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 1.5
$ws.Range("E5").Value = "红色的一把剑"
$ws.Range("F5").Value = "3,4,5"
